$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# --- Row 22: fill in previously-blank payment date/amount (credit now settled) ---
$ws.Range("F22").Value = 44682
$ws.Range("G22").Value = 253

# --- Row 29: fill in previously-blank payment date/amount ---
$ws.Range("F29").Value = 44676
$ws.Range("G29").Value = 21330

# --- Row 30: fill in previously-blank payment date/amount ---
$ws.Range("F30").Value = 44676
$ws.Range("G30").Value = 8587

# --- Row 31: new credit entry ---
$ws.Range("A31").Value = 44676
$ws.Range("D31").Value = "HERRADURA DAVID"
$ws.Range("E31").Value = 11494
$ws.Range("F31").Value = 44679
$ws.Range("G31").Value = 11494

# --- Row 32: new credit entry ---
$ws.Range("A32").Value = 44676
$ws.Range("D32").Value = "OBRADOR"
$ws.Range("E32").Value = 52
$ws.Range("F32").Value = 44682
$ws.Range("G32").Value = 52

# --- Row 33: new credit entry ---
$ws.Range("A33").Value = 44679
$ws.Range("D33").Value = "HERRADURA DAVID"
$ws.Range("E33").Value = 15177
$ws.Range("F33").Value = 44680
$ws.Range("G33").Value = 15177

# --- Row 34: new credit entry ---
$ws.Range("A34").Value = 44679
$ws.Range("D34").Value = "OBRADOR"
$ws.Range("E34").Value = 528
$ws.Range("F34").Value = 44682
$ws.Range("G34").Value = 528

# --- Row 35: new credit entry ---
$ws.Range("A35").Value = 44680
$ws.Range("D35").Value = "EL PRIMO"
$ws.Range("E35").Value = 23340
$ws.Range("F35").Value = 44682
$ws.Range("G35").Value = 23340

# --- Row 36: new credit entry ---
$ws.Range("A36").Value = 44680
$ws.Range("D36").Value = "HERRADURA DAVID"
$ws.Range("E36").Value = 15388
$ws.Range("F36").Value = 44681
$ws.Range("G36").Value = 15388

# --- Row 37: new credit entry ---
$ws.Range("A37").Value = 44681
$ws.Range("D37").Value = "MAURO"
$ws.Range("E37").Value = 6422
$ws.Range("F37").Value = 44682
$ws.Range("G37").Value = 6422

# --- Row 38: new credit entry, still unpaid -> highlight F38:G38 in blue-on-yellow ---
$ws.Range("A38").Value = 44681
$ws.Range("D38").Value = "HERRADURA DAVID"
$ws.Range("E38").Value = 11261
$ws.Range("F38:G38").Font.Color = 16711680
$ws.Range("F38:G38").Interior.Color = 65535

# --- Move the active cell selection to reflect where the user left off editing ---
$ws.Range("D40").Select()
